$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 433. This shifts the existing rows 433-457
# down to 436-460 (preserving all their data/styles).
$ws.Rows("433:435").Insert()

# Fill the 3 newly inserted rows with data for 2022-02-18 (Excel serial 44610).
$data = @(
    @{ Row = 433; I = "Extra";   J = 110; K = 3000; L = 3000; M = 3000; N = "`$/unidad"; O = "Paine"; P = 3000 },
    @{ Row = 434; I = "Primera"; J = 120; K = 2000; L = 2000; M = 2000; N = "`$/unidad"; O = "Paine"; P = 2000 },
    @{ Row = 435; I = "Segunda"; J = 120; K = 1500; L = 1500; M = 1500; N = "`$/unidad"; O = "Paine"; P = 1500 }
)

foreach ($d in $data) {
    $r = $d.Row
    $ws.Cells.Item($r, 1).Value = 3
    $ws.Cells.Item($r, 2).Value = "Femacal de La Calera"
    $ws.Cells.Item($r, 3).Value = "Coquimbo"
    $ws.Cells.Item($r, 4).Value = 44610
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 5).Value = 5
    $ws.Cells.Item($r, 6).Value = 100112028
    $ws.Cells.Item($r, 7).Value = "Sandia"
    $ws.Cells.Item($r, 8).Value = "Sin especificar"
    $ws.Cells.Item($r, 9).Value = $d.I
    $ws.Cells.Item($r, 10).Value = $d.J
    $ws.Cells.Item($r, 11).Value = $d.K
    $ws.Cells.Item($r, 12).Value = $d.L
    $ws.Cells.Item($r, 13).Value = $d.M
    $ws.Cells.Item($r, 14).Value = $d.N
    $ws.Cells.Item($r, 15).Value = $d.O
    $ws.Cells.Item($r, 16).Value = $d.P
    $ws.Cells.Item($r, 17).Value = 1
    $ws.Cells.Item($r, 18).Value = "Hortaliza"
}
